# Lab1.xlsx regression fix
# - "Práctica 1" is Worksheets.Item(1), "Práctica 2" is Worksheets.Item(2)
#
# 1) Formula fixes on "Práctica 2" (sheet 2):
#      E12:E16  use of the non-Excel function pow(x,2) is replaced by (x^2)
#      E22:E25  relative reference to D22 becomes an absolute $D$22
#      E32:E36  now reference their own row's D (col D row 32, absolute $D$32)
#                and use the full PI() angle instead of (2/3)*PI()
# 2) View-state fixes: workbook now opens on "Práctica 2" (which becomes the
#    selected/active tab), with its own scroll position and selection square,
#    while "Práctica 1" is no longer the selected tab (selection itself is
#    unchanged there).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Formula corrections on "Práctica 2"
# ---------------------------------------------------------------------------

# E12:E16 -- replace pow(Cxx, 2) with (Cxx^2)
$ws2.Range("E12").Formula = '=1000*(C12^2)/$D$12*(1-COS(PI()/2))'
$ws2.Range("E13").Formula = '=1000*(C13^2)/$D$12*(1-COS(PI()/2))'
$ws2.Range("E14").Formula = '=1000*(C14^2)/$D$12*(1-COS(PI()/2))'
$ws2.Range("E15").Formula = '=1000*(C15^2)/$D$12*(1-COS(PI()/2))'
$ws2.Range("E16").Formula = '=1000*(C16^2)/$D$12*(1-COS(PI()/2))'

# E22:E25 -- make the reference to D22 absolute
$ws2.Range("E22").Formula = '=1000*((C22^2)/$D$22)*(1-COS((2/3)*PI()))'
$ws2.Range("E23").Formula = '=1000*((C23^2)/$D$22)*(1-COS((2/3)*PI()))'
$ws2.Range("E24").Formula = '=1000*((C24^2)/$D$22)*(1-COS((2/3)*PI()))'
$ws2.Range("E25").Formula = '=1000*((C25^2)/$D$22)*(1-COS((2/3)*PI()))'

# E32:E36 -- reference own-row absolute $D$32 and use full PI() angle
$ws2.Range("E32").Formula = '=1000*((C32^2)/$D$32)*(1-COS((PI())))'
$ws2.Range("E33").Formula = '=1000*((C33^2)/$D$32)*(1-COS((PI())))'
$ws2.Range("E34").Formula = '=1000*((C34^2)/$D$32)*(1-COS((PI())))'
$ws2.Range("E35").Formula = '=1000*((C35^2)/$D$32)*(1-COS((PI())))'
$ws2.Range("E36").Formula = '=1000*((C36^2)/$D$32)*(1-COS((PI())))'

# Row heights tighten slightly on every row touched by the recalculation
# above (12-16, 22-26, 32-36) once the workbook re-lays itself out.
$heightRows = @(12,13,14,15,16,22,23,24,25,26,32,33,34,35,36)
foreach ($r in $heightRows) {
    $ws2.Rows.Item($r).RowHeight = 13.8
}

# ---------------------------------------------------------------------------
# 2) View / selection state: "Práctica 2" becomes the active, selected sheet
# ---------------------------------------------------------------------------

$ws1.Range("P5").Select()
$ws2.Activate()
$ws2.Range("C32").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
